# Zafina Tekken 8 frame data: invert the sign of the "Block" (column E)
# values for all data rows (2..137). Rows whose Block cell is blank are
# skipped automatically since there is no numeric value to negate.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 137; $row++) {
    $cell = $ws.Cells.Item($row, 5)   # Column E = Block
    $val = $cell.Value2
    if ($val -ne $null -and $val -is [double]) {
        $cell.Value = -$val
    }
}
